$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.800.85'
$ws.Range('E2').Value = '  +4.05%  '
$ws.Range('D3').Value = '2.657.55'
$ws.Range('E3').Value = '  +1.37%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '''569.88'
$ws.Range('E5').Value = '  +6.84%  '
$ws.Range('D6').Value = '''148.20'
$ws.Range('E6').Value = '  +4.13%  '
$ws.Range('E7').Value = '  -0.48%  '
$ws.Range('D8').Value = '''0.607'
$ws.Range('E8').Value = '  +7.06%  '
$ws.Range('D9').Value = '''6.86'
$ws.Range('E9').Value = '  -1.07%  '
$ws.Range('E10').Value = '  +4.82%  '
$ws.Range('E11').Value = '  +6.65%  '
$ws.Range('D12').Value = '''0.345'
$ws.Range('E12').Value = '  +3.48%  '
$ws.Range('D13').Value = '3.128.18'
$ws.Range('E13').Value = '  +1.26%  '
$ws.Range('D14').Value = '60.759.32'
$ws.Range('E14').Value = '  +4.09%  '
$ws.Range('D15').Value = '''21.94'
$ws.Range('E15').Value = '  +5.97%  '
$ws.Range('D16').Value = '''0.0000138'
$ws.Range('E16').Value = '  +4.86%  '
$ws.Range('D17').Value = '2.666.67'
$ws.Range('E17').Value = '  +1.32%  '
$ws.Range('D18').Value = '''4.57'
$ws.Range('E18').Value = '  +4.11%  '
$ws.Range('D19').Value = '''346.44'
$ws.Range('E19').Value = '  +3.71%  '
$ws.Range('D20').Value = '''10.50'
$ws.Range('E20').Value = '  +3.75%  '
$ws.Range('D21').Value = '''6.46'
$ws.Range('E21').Value = '  +3.87%  '
$ws.Range('D22').Value = '''5.84'
$ws.Range('D23').Value = '''0.997'
$ws.Range('E23').Value = '  -0.12%  '
$ws.Range('D24').Value = '''67.01'
$ws.Range('E24').Value = '  +1.07%  '
$ws.Range('D25').Value = '''0.443'
$ws.Range('E25').Value = '  +6.75%  '
$ws.Range('D26').Value = '''0.167'
$ws.Range('E26').Value = '  +2.71%  '
$ws.Range('D27').Value = '''0.990'
$ws.Range('E27').Value = '  -0.86%  '
$ws.Range('D28').Value = '''7.41'
$ws.Range('E28').Value = '  +4.59%  '
$ws.Range('D29').Value = '0.0₃0796'
$ws.Range('E29').Value = '  +8.20%  '
$ws.Range('E30').Value = '  -0.20%  '
$ws.Range('E31').Value = '  +5.08%  '
$ws.Range('D32').Value = '''6.16'
$ws.Range('E32').Value = '  +5.41%  '
$ws.Range('B33').Value = 'EthereumClassic'
$ws.Range('C33').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D33').Value = '''19.37'
$ws.Range('E33').Value = '  +3.28%  '
$ws.Range('B34').Value = 'Monero'
$ws.Range('C34').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D34').Value = '''155.50'
$ws.Range('E34').Value = '  +3.49%  '
$ws.Range('D35').Value = '''4.13'
$ws.Range('E35').Value = '  +6.16%  '
$ws.Range('B36').Value = 'SuiNetwork'
$ws.Range('C36').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D36').Value = '''0.926'
$ws.Range('E36').Value = '  +8.98%  '
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').Value = '''1.19'
$ws.Range('E37').Value = '  +8.41%  '
$ws.Range('B38').Value = 'Fetch.AI'
$ws.Range('C38').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D38').Value = '''0.921'
$ws.Range('E38').Value = '  +14.04%  '
$ws.Range('D39').Value = '''37.72'
$ws.Range('E39').Value = '  +1.50%  '
$ws.Range('D40').Value = '''1.53'
$ws.Range('E40').Value = '  +8.43%  '
$ws.Range('D41').Value = '''311.96'
$ws.Range('E41').Value = '  +11.37%  '
$ws.Range('D42').Value = '''3.69'
$ws.Range('E42').Value = '  +3.67%  '
$ws.Range('D43').Value = '''0.612'
$ws.Range('E43').Value = '  +2.87%  '
$ws.Range('E44').Value = '  -0.59%  '
$ws.Range('D45').Value = '''0.0984'
$ws.Range('E45').Value = '  +5.24%  '
$ws.Range('D46').Value = '''0.0553'
$ws.Range('E46').Value = '  +4.78%  '
$ws.Range('D47').Value = '''19.70'
$ws.Range('E47').Value = '  +3.92%  '
$ws.Range('D48').Value = '''10.67'
$ws.Range('E48').Value = '  -0.19%  '
$ws.Range('E49').Value = '  +5.57%  '
$ws.Range('D50').Value = '''125.96'
$ws.Range('E50').Value = '  +11.67%  '
$ws.Range('B51').Value = 'RenderToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D51').Value = '''4.78'
$ws.Range('E51').Value = '  +7.84%  '
